# "+ balance after deleting material"
#
# The material list grew by 3 new rows (materials #7, #8 and #9) and the
# trailing "balance" rows shifted down from rows 8-10 to rows 11-13 with
# recomputed totals. Several existing rows were also relabeled/retotaled.
#
# NOTE: every cell on this sheet is stored as literal text (even the ones
# that look like plain numbers or dates), so any value that Excel would
# otherwise auto-convert to a number/date/currency is entered with a
# leading apostrophe to force text entry, exactly as a user typing into a
# pre-existing text cell would do.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: "fwe" -> "ttttt", price zeroed out ---
$ws.Range("B3").Value = "ttttt"
$ws.Range("E3").Value = "'0"
$ws.Range("F3").Value = "'0.0"
$ws.Range("G3").Value = " "

# --- Row 4: "ttttt" -> "Material1", unit кг -> м, qty/date updated ---
$ws.Range("B4").Value = "Material1"
$ws.Range("C4").Value = "м"
$ws.Range("E4").Value = "'10"
$ws.Range("F4").Value = "'10.0"
$ws.Range("H4").Value = "'01.05.2021"

# --- Row 5: "Material1" -> "material2", unit м -> кг, qty updated ---
$ws.Range("B5").Value = "material2"
$ws.Range("C5").Value = "кг"
$ws.Range("E5").Value = "'50"
$ws.Range("F5").Value = "'50.0"

# --- Row 6: "material2" -> "rfreuwgfpslw", qty/date updated ---
$ws.Range("B6").Value = "rfreuwgfpslw"
$ws.Range("E6").Value = "'5000"
$ws.Range("F6").Value = "'5000.0"
$ws.Range("H6").Value = "'17.05.2021"

# --- Row 7: "material3" -> "dedede", unit м^3 -> кг, qty/price/date updated ---
$ws.Range("B7").Value = "dedede"
$ws.Range("C7").Value = "кг"
$ws.Range("D7").Value = "'1"
$ws.Range("E7").Value = "'5000"
$ws.Range("F7").Value = "'5000.0"
$ws.Range("G7").Value = " "
$ws.Range("H7").Value = "'17.05.2021"

# --- Row 8 (new material #7: "wdwdw") ---
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "wdwdw"
$ws.Range("C8").Value = "м"
$ws.Range("D8").Value = "'1"
$ws.Range("E8").Value = "'5000"
$ws.Range("F8").Value = "'5000.0"
$ws.Range("G8").Value = " "
$ws.Range("H8").Value = "'17.05.2021"

# --- Row 9 (new material #8: "frfewrf") ---
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "frfewrf"
$ws.Range("C9").Value = "м"
$ws.Range("D9").Value = "'1"
$ws.Range("E9").Value = "'1"
$ws.Range("F9").Value = " "
$ws.Range("G9").Value = "'1.0"
$ws.Range("H9").Value = "'17.05.2021"

# --- Row 10 (new material #9: "kjgrguregfiuesas") ---
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "kjgrguregfiuesas"
$ws.Range("C10").Value = "кг"
$ws.Range("D10").Value = "'1"
$ws.Range("E10").Value = "'6"
$ws.Range("F10").Value = " "
$ws.Range("G10").Value = "'6.0"
$ws.Range("H10").Value = "'18.05.2021"

# --- Row 11 (blank + recomputed sum/($) totals, was row 8) ---
$ws.Range("A11").Value = " "
$ws.Range("B11").Value = " "
$ws.Range("C11").Value = " "
$ws.Range("D11").Value = " "
$ws.Range("E11").Value = " "
$ws.Range("F11").Value = "'15588.0"
$ws.Range("G11").Value = "'7.0"
$ws.Range("H11").Value = " "

# --- Row 12 (blank + recomputed exchange-rate row, was row 9) ---
$ws.Range("A12").Value = " "
$ws.Range("B12").Value = " "
$ws.Range("C12").Value = " "
$ws.Range("D12").Value = " "
$ws.Range("E12").Value = "'10000.0"
$ws.Range("F12").Value = "'1.5588"
$ws.Range("G12").Value = " "
$ws.Range("H12").Value = " "

# --- Row 13 (blank + recomputed "$" balance row, was row 10) ---
$ws.Range("A13").Value = " "
$ws.Range("B13").Value = " "
$ws.Range("C13").Value = " "
$ws.Range("D13").Value = " "
$ws.Range("E13").Value = " "
$ws.Range("F13").Value = "'`$8.5588"
$ws.Range("G13").Value = " "
$ws.Range("H13").Value = " "

# Give the new "#" cells (A11:A13) the same bold/centered/thin-border style
# already used by the "#" column's blank summary rows (A8:A10 originally).
$ws.Range("A10").Copy()
$ws.Range("A11:A13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
